$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the sheet name (remove the "(2010-2023)" suffix)
$ws.Name = "g3.5c Aumento"

# Add the new "Ano" column header with the same style as the other headers
$ws.Range("D1").Value = "Ano"
$ws.Range("D1").Style = $ws.Range("A1").Style

# Fill the "Ano" column with the year range for each data row
$ws.Range("D2:D9").Value = "2010-2023"
